{"js": "// Update the worksheet date line and regenerate all 100 addition/subtraction\n// answers in the 20x5 practice table, per commit \"Update master to output\n// generated at 9a8706d\".\n\n// 1) Update the date paragraph (\"2024-01-12 Friday\" -> \"2024-01-13 Saturday\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\nif (dateParagraph.text === \"2024-01-12 Friday\") {\n  dateParagraph.insertText(\"2024-01-13 Saturday\", \"Replace\");\n}\n\n// 2) Update every cell of the first (only) table with its new expression,\n// addressed by row/column so the two duplicated \"before\" strings\n// (\"2+22=24\" and \"39+27=66\") each resolve to the correct distinct \"after\"\n// value instead of a blanket find/replace.\nconst newValues = [\n  [\"56-12=44\", \"84-78=6\", \"69-17=52\", \"1+9=10\", \"87-11=76\"],\n  [\"82-35=47\", \"99-14=85\", \"94-56=38\", \"22+31=53\", \"94-55=39\"],\n  [\"57-36=21\", \"97-3=94\", \"10-2=8\", \"16+83=99\", \"9+42=51\"],\n  [\"70+25=95\", \"92-38=54\", \"26-12=14\", \"97-71=26\", \"33+42=75\"],\n  [\"80-5=75\", \"16+13=29\", \"15+17=32\", \"53-35=18\", \"65+16=81\"],\n  [\"89-40=49\", \"22-15=7\", \"38-25=13\", \"67+16=83\", \"62-22=40\"],\n  [\"82-5=77\", \"24+20=44\", \"42-0=42\", \"84-79=5\", \"89-63=26\"],\n  [\"44+7=51\", \"7+59=66\", \"21-12=9\", \"1+55=56\", \"57+2=59\"],\n  [\"10+22=32\", \"95+0=95\", \"25+46=71\", \"97-80=17\", \"65+25=90\"],\n  [\"91-85=6\", \"29-20=9\", \"74-34=40\", \"95-51=44\", \"71-63=8\"],\n  [\"40+34=74\", \"79+17=96\", \"4+72=76\", \"18+13=31\", \"66-0=66\"],\n  [\"72-26=46\", \"61+12=73\", \"91-44=47\", \"0+53=53\", \"70-13=57\"],\n  [\"69-52=17\", \"59-42=17\", \"49+0=49\", \"14-9=5\", \"27+68=95\"],\n  [\"92-77=15\", \"12+80=92\", \"19+80=99\", \"39+25=64\", \"77+14=91\"],\n  [\"86-11=75\", \"86-51=35\", \"35+51=86\", \"38+56=94\", \"5+15=20\"],\n  [\"3+20=23\", \"83-26=57\", \"10+24=34\", \"99-10=89\", \"72-5=67\"],\n  [\"86-16=70\", \"57+10=67\", \"84-77=7\", \"79-11=68\", \"65+2=67\"],\n  [\"52+30=82\", \"37+2=39\", \"93-20=73\", \"36-33=3\", \"96-33=63\"],\n  [\"12+63=75\", \"76-75=1\", \"71-26=45\", \"2+97=99\", \"0+53=53\"],\n  [\"96-0=96\", \"62+24=86\", \"61-44=17\", \"58-37=21\", \"27+18=45\"],\n];\n\nconst table = body.tables.getFirst();\ntable.values = newValues;\n\nawait context.sync();\n", "ps1": "# Update master to output generated at 9a8706d\n# Updates the worksheet date line and regenerates all 100 addition/\n# subtraction answers in the 20x5 practice table.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph.\n$dateParagraph = $d.Paragraphs.Item(1)\nif ($dateParagraph.Range.Text.TrimEnd([char]13) -eq \"2024-01-12 Friday\") {\n  $dateParagraph.Range.Text = \"2024-01-13 Saturday\"\n}\n\n# 2) Update every cell of the first (only) table with its new expression,\n#    addressed by row/column (1-based, like Word COM) so the two\n#    duplicated \"before\" strings (\"2+22=24\" and \"39+27=66\") each resolve\n#    to the correct distinct \"after\" value.\n$newValues = @(\n  @(\"56-12=44\", \"84-78=6\", \"69-17=52\", \"1+9=10\", \"87-11=76\"),\n  @(\"82-35=47\", \"99-14=85\", \"94-56=38\", \"22+31=53\", \"94-55=39\"),\n  @(\"57-36=21\", \"97-3=94\", \"10-2=8\", \"16+83=99\", \"9+42=51\"),\n  @(\"70+25=95\", \"92-38=54\", \"26-12=14\", \"97-71=26\", \"33+42=75\"),\n  @(\"80-5=75\", \"16+13=29\", \"15+17=32\", \"53-35=18\", \"65+16=81\"),\n  @(\"89-40=49\", \"22-15=7\", \"38-25=13\", \"67+16=83\", \"62-22=40\"),\n  @(\"82-5=77\", \"24+20=44\", \"42-0=42\", \"84-79=5\", \"89-63=26\"),\n  @(\"44+7=51\", \"7+59=66\", \"21-12=9\", \"1+55=56\", \"57+2=59\"),\n  @(\"10+22=32\", \"95+0=95\", \"25+46=71\", \"97-80=17\", \"65+25=90\"),\n  @(\"91-85=6\", \"29-20=9\", \"74-34=40\", \"95-51=44\", \"71-63=8\"),\n  @(\"40+34=74\", \"79+17=96\", \"4+72=76\", \"18+13=31\", \"66-0=66\"),\n  @(\"72-26=46\", \"61+12=73\", \"91-44=47\", \"0+53=53\", \"70-13=57\"),\n  @(\"69-52=17\", \"59-42=17\", \"49+0=49\", \"14-9=5\", \"27+68=95\"),\n  @(\"92-77=15\", \"12+80=92\", \"19+80=99\", \"39+25=64\", \"77+14=91\"),\n  @(\"86-11=75\", \"86-51=35\", \"35+51=86\", \"38+56=94\", \"5+15=20\"),\n  @(\"3+20=23\", \"83-26=57\", \"10+24=34\", \"99-10=89\", \"72-5=67\"),\n  @(\"86-16=70\", \"57+10=67\", \"84-77=7\", \"79-11=68\", \"65+2=67\"),\n  @(\"52+30=82\", \"37+2=39\", \"93-20=73\", \"36-33=3\", \"96-33=63\"),\n  @(\"12+63=75\", \"76-75=1\", \"71-26=45\", \"2+97=99\", \"0+53=53\"),\n  @(\"96-0=96\", \"62+24=86\", \"61-44=17\", \"58-37=21\", \"27+18=45\"),\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n  for ($c = 0; $c -lt $newValues[$r].Count; $c++) {\n    $t.Cell($r + 1, $c + 1).Range.Text = $newValues[$r][$c]\n  }\n}\n\n"}
